$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Fill in the description for the YAX3MAUTHLIB row (B9), which was previously
# left blank and highlighted green while still being worked on.
$ws.Range("B9").Value = "Libreria per la gestione dell'Autenticazione"

# The row is no longer a "work in progress" row, so drop the green highlight
# fill that had been applied to it. Borrow the (unfilled, bordered) formatting
# from the neighboring cell in the same column so the border styling is kept.
$ws.Range("C9").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# A15 also still carried a stray "no fill" override left over from earlier
# editing; normalize it to match the same plain formatting used by the rows
# around it.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the active selection to C9 (matches the sheetView selection change).
$ws.Range("C9").Select()
